$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$lb = [char]11

$t.Cell(1,1).Range.Text = "82 x 39" + $lb + "  3    9" + $lb + "  ----" + $lb + "8|    |" + $lb + "2|    |"
$t.Cell(1,2).Range.Text = "12 x 86" + $lb + "  8    6" + $lb + "  ----" + $lb + "1|    |" + $lb + "2|    |"
$t.Cell(1,3).Range.Text = "64 x 73" + $lb + "  7    3" + $lb + "  ----" + $lb + "6|    |" + $lb + "4|    |"
$t.Cell(2,1).Range.Text = "60 x 63" + $lb + "  6    3" + $lb + "  ----" + $lb + "6|    |" + $lb + "0|    |"
$t.Cell(2,2).Range.Text = "27 x 98" + $lb + "  9    8" + $lb + "  ----" + $lb + "2|    |" + $lb + "7|    |"
$t.Cell(2,3).Range.Text = "18 x 54" + $lb + "  5    4" + $lb + "  ----" + $lb + "1|    |" + $lb + "8|    |"
$t.Cell(3,1).Range.Text = "10 x 76" + $lb + "  7    6" + $lb + "  ----" + $lb + "1|    |" + $lb + "0|    |"
$t.Cell(3,2).Range.Text = "54 x 21" + $lb + "  2    1" + $lb + "  ----" + $lb + "5|    |" + $lb + "4|    |"
$t.Cell(3,3).Range.Text = "94 x 83" + $lb + "  8    3" + $lb + "  ----" + $lb + "9|    |" + $lb + "4|    |"
$t.Cell(4,1).Range.Text = "52 x 79" + $lb + "  7    9" + $lb + "  ----" + $lb + "5|    |" + $lb + "2|    |"
$t.Cell(4,2).Range.Text = "63 x 91" + $lb + "  9    1" + $lb + "  ----" + $lb + "6|    |" + $lb + "3|    |"
$t.Cell(4,3).Range.Text = "40 x 79" + $lb + "  7    9" + $lb + "  ----" + $lb + "4|    |" + $lb + "0|    |"
$t.Cell(5,1).Range.Text = "32 x 83" + $lb + "  8    3" + $lb + "  ----" + $lb + "3|    |" + $lb + "2|    |"
$t.Cell(5,2).Range.Text = "20 x 26" + $lb + "  2    6" + $lb + "  ----" + $lb + "2|    |" + $lb + "0|    |"
$t.Cell(5,3).Range.Text = "63 x 22" + $lb + "  2    2" + $lb + "  ----" + $lb + "6|    |" + $lb + "3|    |"
